# Update market/profit data across multiple sheets (scheduled runner refresh).
# Applies updated currentAveragePrice / LevePrice / LeveProfit values per the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 2045.2157
$ws.Range("I38").Value = 359.7037
$ws.Range("J38").Value = 3941.4167
$ws.Range("K38").Value = 1079.1111
$ws.Range("L38").Value = 11824.2501
$ws.Range("M38").Value = -707.1111000000001
$ws.Range("N38").Value = -12568.2501

$ws.Range("H43").Value = 2388.3333
$ws.Range("I43").Value = 2400
$ws.Range("J43").Value = 2386.875
$ws.Range("K43").Value = 2400
$ws.Range("L43").Value = 2386.875
$ws.Range("M43").Value = -2331
$ws.Range("N43").Value = -2524.875

$ws.Range("H86").Value = 1462.125
$ws.Range("I86").Value = 1445.7273
$ws.Range("J86").Value = 1498.2
$ws.Range("K86").Value = 1445.7273
$ws.Range("L86").Value = 1498.2
$ws.Range("M86").Value = -322.7273
$ws.Range("N86").Value = -3744.2

$ws.Range("H89").Value = 1462.125
$ws.Range("I89").Value = 1445.7273
$ws.Range("J89").Value = 1498.2
$ws.Range("K89").Value = 7228.636500000001
$ws.Range("L89").Value = 7491
$ws.Range("M89").Value = -1612.636500000001
$ws.Range("N89").Value = -18723

$ws.Range("H118").Value = 1350.1666
$ws.Range("I118").Value = 325.25
$ws.Range("J118").Value = 3400
$ws.Range("K118").Value = 975.75
$ws.Range("L118").Value = 10200
$ws.Range("M118").Value = 681.25
$ws.Range("N118").Value = -13514

$ws.Range("H137").Value = 1742.5454
$ws.Range("I137").Value = 1302.6428
$ws.Range("J137").Value = 2512.375
$ws.Range("K137").Value = 3907.9284
$ws.Range("L137").Value = 7537.125
$ws.Range("M137").Value = -1357.9284
$ws.Range("N137").Value = -12637.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1254599.1
$ws.Range("I32").Value = 12445.788
$ws.Range("J32").Value = 21750130
$ws.Range("K32").Value = 12445.788
$ws.Range("L32").Value = 21750130
$ws.Range("M32").Value = -12158.788
$ws.Range("N32").Value = -21750704

$ws.Range("H37").Value = 6303.5
$ws.Range("J37").Value = 8238
$ws.Range("L37").Value = 8238
$ws.Range("N37").Value = -8784

$ws.Range("H45").Value = 2646.6765
$ws.Range("I45").Value = 2379.926
$ws.Range("J45").Value = 3675.5715
$ws.Range("K45").Value = 2379.926
$ws.Range("L45").Value = 3675.5715
$ws.Range("M45").Value = -2002.926
$ws.Range("N45").Value = -4429.5715

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1986.6333
$ws.Range("I105").Value = 2029.5927
$ws.Range("J105").Value = 1600
$ws.Range("K105").Value = 2029.5927
$ws.Range("L105").Value = 1600
$ws.Range("M105").Value = -282.5926999999999
$ws.Range("N105").Value = -5094

$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 12484.333
$ws.Range("J59").Value = 12484.333
$ws.Range("L59").Value = 12484.333
$ws.Range("N59").Value = -14774.333

$ws.Range("H60").Value = 7033.8887
$ws.Range("J60").Value = 8262.4
$ws.Range("L60").Value = 8262.4
$ws.Range("N60").Value = -9284.4

$ws.Range("H68").Value = 15492.714
$ws.Range("J68").Value = 15492.714
$ws.Range("L68").Value = 15492.714
$ws.Range("N68").Value = -16990.714

$ws.Range("H71").Value = 15492.714
$ws.Range("J71").Value = 15492.714
$ws.Range("L71").Value = 46478.142
$ws.Range("N71").Value = -53966.142

$ws.Range("H74").Value = 17928.334
$ws.Range("J74").Value = 17928.334
$ws.Range("L74").Value = 17928.334
$ws.Range("N74").Value = -19676.334

$ws.Range("H77").Value = 17928.334
$ws.Range("J77").Value = 17928.334
$ws.Range("L77").Value = 53785.00199999999
$ws.Range("N77").Value = -62521.00199999999

$ws.Range("H105").Value = 1363
$ws.Range("I105").Value = 1363
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1363
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 384
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 344.08694
$ws.Range("I23").Value = 226.07143
$ws.Range("J23").Value = 527.6667
$ws.Range("K23").Value = 678.21429
$ws.Range("L23").Value = 1583.0001
$ws.Range("M23").Value = -443.21429
$ws.Range("N23").Value = -2053.0001

$ws.Range("H68").Value = 1024.7222
$ws.Range("I68").Value = 828.8570999999999
$ws.Range("J68").Value = 1149.3636
$ws.Range("K68").Value = 2486.5713
$ws.Range("L68").Value = 3448.0908
$ws.Range("M68").Value = -1675.5713
$ws.Range("N68").Value = -5070.0908

$ws.Range("H71").Value = 1024.7222
$ws.Range("I71").Value = 828.8570999999999
$ws.Range("J71").Value = 1149.3636
$ws.Range("K71").Value = 7459.7139
$ws.Range("L71").Value = 10344.2724
$ws.Range("M71").Value = -3403.7139
$ws.Range("N71").Value = -18456.2724

$ws.Range("H122").Value = 594.0294
$ws.Range("I122").Value = 418.25925
$ws.Range("J122").Value = 1272
$ws.Range("K122").Value = 3764.33325
$ws.Range("L122").Value = 11448
$ws.Range("M122").Value = -1314.33325
$ws.Range("N122").Value = -16348

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1609.4688
$ws.Range("I97").Value = 1001.13336
$ws.Range("J97").Value = 2146.2354
$ws.Range("K97").Value = 1001.13336
$ws.Range("L97").Value = 2146.2354
$ws.Range("M97").Value = -505.13336
$ws.Range("N97").Value = -3138.2354

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6409.091
$ws.Range("J62").Value = 6333.6665
$ws.Range("L62").Value = 6333.6665
$ws.Range("N62").Value = -7581.6665

$ws.Range("H65").Value = 6409.091
$ws.Range("J65").Value = 6333.6665
$ws.Range("L65").Value = 31668.3325
$ws.Range("N65").Value = -37908.3325

$ws.Range("H81").Value = 66669936
$ws.Range("I81").Value = 2851.077
$ws.Range("J81").Value = 500006000
$ws.Range("K81").Value = 5702.154
$ws.Range("L81").Value = 1000012000
$ws.Range("M81").Value = -4641.154
$ws.Range("N81").Value = -1000014122

$ws.Range("H84").Value = 66669936
$ws.Range("I84").Value = 2851.077
$ws.Range("J84").Value = 500006000
$ws.Range("K84").Value = 28510.77
$ws.Range("L84").Value = 5000060000
$ws.Range("M84").Value = -23206.77
$ws.Range("N84").Value = -5000070608

